# data drive test case 1
# Add a new "ValidLogin" worksheet after the existing "TC1" sheet and
# populate it with username/password test data, then make it the active
# (selected) sheet/tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after TC1.
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "ValidLogin"

# Data-drive the login test case.
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "pointofsale"

# Make ValidLogin the active sheet/tab, zoomed in, with B3 selected.
$ws2.Activate()
$excel.ActiveWindow.Zoom = 160
[void]$ws2.Range("B3").Select()
